$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells for new columns
$ws.Range("AC1").Value = "Wins"
$ws.Range("AD1").Value = "Losses"
$ws.Range("AE1").Value = "Ties"

# Copy header style (bold, border, centered) from an existing header cell (e.g. AB1)
$ws.Range("AB1").Copy()
$ws.Range("AC1:AE1").PasteSpecial(-4122) # xlPasteFormats

# Fill in the Wins/Losses/Ties values for all data rows (2-46)
for ($r = 2; $r -le 46; $r++) {
    $ws.Cells.Item($r, 29).Value = 72  # AC
    $ws.Cells.Item($r, 30).Value = 90  # AD
    $ws.Cells.Item($r, 31).Value = 0   # AE
}
